# switched LS and HS outputs to match hardware timers distribution.
# Added a 2k pulldown for the tach output -> new resistor R17 joins the
# existing 2K pulldown designator group in the BOM sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "R17,R16,R15,R14,R13,R12,R11,R10,R9,R8,R7,R6,R5,R4,R3,R2,R1"

$ws.Range("B3").Select()
